$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TSCO")

# Insert a new column before column D, shifting existing D:K data to E:L.
$ws.Range("D1").EntireColumn.Insert()

# Give the freshly-inserted column D the same formatting as column E (which
# holds the data that used to live in D), so number formats/fonts line up.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)

# New (2018) financial-year data added into the freshly inserted column D.
$newData = @{
    7  = 43463
    8  = 7911000
    9  = 5208500
    10 = 2702500
    12 = "NA"
    13 = 0
    14 = 0
    15 = 177400
    17 = 7209300
    18 = 701700
    20 = 100
    21 = 879100
    22 = 18400
    23 = 683400
    24 = 151000
    25 = 0
    26 = 532400
    27 = 532400
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = -100
    33 = 532400
    34 = 0
    35 = 532400
    38 = 43463
    41 = 86300
    42 = 0
    43 = 4100
    44 = 1589500
    45 = 114400
    46 = 1794400
    47 = 0
    48 = 1134500
    49 = 124500
    50 = 0
    51 = 0
    52 = 31900
    53 = 0
    54 = 3085300
    57 = 620000
    58 = 29900
    59 = 288200
    60 = 938100
    61 = 410400
    62 = 175000
    63 = 0
    64 = 0
    65 = 0
    66 = 1523400
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = 3213900
    73 = 0
    74 = 0
    75 = 0
    76 = 1561800
    77 = 0
    80 = 43463
    81 = 532400
    83 = 177400
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = 694400
    91 = -278500
    92 = 0
    93 = 0
    94 = -276300
    96 = -147100
    97 = 0
    98 = 0
    99 = 0
    100 = -440900
    101 = 0
    102 = -22800
}

foreach ($row in $newData.Keys) {
    $ws.Cells.Item($row, 4).Value = $newData[$row]
}
